$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Price/Volume columns so Excel does not
# auto-convert numeric-looking strings (e.g. "1.00") into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "36.399.55"
$ws.Range("E2").Value = "  -2.23%  "
$ws.Range("D3").Value = "2.044.70"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "241.29"
$ws.Range("E5").Value = "  -3.26%  "
$ws.Range("D6").Value = "0.665"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "54.47"
$ws.Range("E8").Value = "  -7.48%  "
$ws.Range("D9").Value = "58.27"
$ws.Range("E9").Value = "  -3.15%  "
$ws.Range("D10").Value = "0.355"
$ws.Range("E10").Value = "  -8.09%  "
$ws.Range("D11").Value = "0.0744"
$ws.Range("E11").Value = "  -5.93%  "
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.884"
$ws.Range("E13").Value = "  -3.44%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "14.66"
$ws.Range("E14").Value = "  -8.89%  "
$ws.Range("D15").Value = "2.343.24"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "5.33"
$ws.Range("E16").Value = "  -7.98%  "
$ws.Range("D17").Value = "2.049.24"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "36.365.25"
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("D19").Value = "16.45"
$ws.Range("E19").Value = "  -14.10%  "
$ws.Range("D20").Value = "71.71"
$ws.Range("E20").Value = "  -5.32%  "
$ws.Range("D21").Value = "0.0₃0847"
$ws.Range("E21").Value = "  -6.64%  "
$ws.Range("D22").Value = "236.63"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "5.21"
$ws.Range("E23").Value = "  -5.43%  "
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  -5.85%  "
$ws.Range("D26").Value = "9.25"
$ws.Range("E26").Value = "  -3.63%  "
$ws.Range("D27").Value = "2.10"
$ws.Range("E27").Value = "  -5.32%  "
$ws.Range("D28").Value = "162.20"
$ws.Range("E28").Value = "  -5.50%  "
$ws.Range("D29").Value = "19.98"
$ws.Range("E29").Value = "  -1.83%  "
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("D31").Value = "5.06"
$ws.Range("E31").Value = "  -7.43%  "
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").Value = "4.47"
$ws.Range("E33").Value = "  -5.47%  "
$ws.Range("D34").Value = "0.0589"
$ws.Range("E34").Value = "  -6.67%  "
$ws.Range("D36").Value = "1.86"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.0829"
$ws.Range("E37").Value = "  -5.79%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.15"
$ws.Range("E38").Value = "  -7.99%  "
$ws.Range("D39").Value = "1.24"
$ws.Range("E39").Value = "  -8.02%  "
$ws.Range("D40").Value = "4.77"
$ws.Range("E40").Value = "  -7.66%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0213"
$ws.Range("E41").Value = "  -6.40%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "1.10"
$ws.Range("E42").Value = "  -6.41%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "2.78"
$ws.Range("E43").Value = "  -11.65%  "
$ws.Range("D44").Value = "92.80"
$ws.Range("E44").Value = "  -8.97%  "
$ws.Range("D45").Value = "0.0892"
$ws.Range("E45").Value = "  -12.93%  "
$ws.Range("D46").Value = "1.374.77"
$ws.Range("E46").Value = "  +4.60%  "
$ws.Range("D47").Value = "15.49"
$ws.Range("E47").Value = "  -10.50%  "
$ws.Range("D48").Value = "7.24"
$ws.Range("E48").Value = "  +4.35%  "
$ws.Range("D49").Value = "2.82"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("E50").Value = "  -8.06%  "
$ws.Range("D51").Value = "2.231.46"
$ws.Range("E51").Value = "  -1.03%  "
